# Average function complexities has been added to xlsxContracts data.
# This inserts a new "TOTAL" placeholder row right after the header row
# (row 2), shifting all existing data rows down by one, and removes the
# old trailing TOTAL row (whose real totals are no longer needed at the
# bottom of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 2 (pushes all data down by one row).
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the placeholder TOTAL values.
$ws.Range("A2").Value = "TOTAL"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# The original TOTAL row (formerly row 72) has now shifted down to row 73.
# Remove it so the sheet ends with the last data row (row 72).
$ws.Rows.Item(73).Delete()
